# Pinout.xlsx update:
#  - Replaced the "WATER (PORTB 0~2)" sensor block with a single
#    "WATER_SENS_ANALOG" pin (and its pin range becomes 0~2), which makes the
#    old "WATER (PORTG 0~2)" / "39-41" backup-Arduino water row (row 10)
#    redundant, so that row is removed entirely.
#  - Column A is widened to fit its (now differently-sized) longest label.
#  - Selection cursor left on G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Backup Arduino block (right-hand side, row 4): swap the water-sensor
# label for the new analog sensor pin, and update its pin range.
$ws.Range("E4").Value = "WATER_SENS_ANALOG"
$ws.Range("F4").Value = "0~2"

# The dedicated "WATER (PORTG 0~2)" row for the backup Arduino is no longer
# needed, drop it.
$ws.Range("A10:B10").ClearContents()

# Column A now best-fits its content after the edits above.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where the author left it after editing.
$ws.Range("G11").Select()
